$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1,  "PONTIANAK",     "PP", 2700000, 2700000,  860000,  860000),
    @(2,  "KETAPANG",      "PP", 4900000, 4900000, 3260000, 3260000),
    @(3,  "KAYONG UTARA",  "PP", 5100000, 5100000, 3460000, 3460000),
    @(4,  "SINTANG",       "PP",  300000,  300000,  300000,  300000),
    @(5,  "SEKADAU",       "PP",  400000,  400000,  400000,  400000),
    @(6,  "SANGGAU",       "PP",  500000,  500000,  500000,  500000),
    @(7,  "LANDAK",        "PP",  600000,  600000,  600000,  600000),
    @(8,  "BENGKAYANG",    "PP", 1100000, 1100000, 1100000, 1100000),
    @(9,  "SAMBAS",        "PP", 1200000, 1200000, 1200000, 1200000),
    @(10, "SINGKAWANG",    "PP", 1000000, 1000000, 1000000, 1000000),
    @(11, "MEMPAWAH",      "PP",  800000,  800000,  800000,  800000),
    @(12, "KAPUAS HULU",   "PP",  700000,  700000,  700000,  700000),
    @(13, "KUBU RAYA",     "PP", 2700000, 2700000,  860000,  860000)
)

$row = 3
foreach ($d in $data) {
    $ws.Cells.Item($row, 1).Value = $d[0]
    $ws.Cells.Item($row, 2).Value = $d[1]
    $ws.Cells.Item($row, 3).Value = $d[2]
    $ws.Cells.Item($row, 4).Value = $d[3]
    $ws.Cells.Item($row, 5).Value = $d[4]
    $ws.Cells.Item($row, 6).Value = $d[5]
    $ws.Cells.Item($row, 7).Value = $d[6]
    $row++
}
